# Add Q_7 / Q_8 rows to the SIQ sheet (scope & functional requirements),
# matching the formatting applied by hand in Excel (bigger font sizes on
# the newly typed answers, Arial for the "reserve button" question/answer,
# and the workbook's theme "minor" font - Aptos narrow - for the second
# new question).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Q_7 ---------------------------------------------------------
$ws.Range("A8").Value = "Q_7"

$ws.Range("B8").Value = 'Car reservation by "reserve" button only?'
$ws.Range("B8").Font.Name = "Arial"
$ws.Range("B8").Font.Size = 20

$ws.Range("C8").Value = "yes"
$ws.Range("C8").Font.Name = "Arial"

# --- Row 9: Q_8 ---------------------------------------------------------
$ws.Range("A9").Value = "Q_8"

$ws.Range("B9").Value = "Car can be reserved by one person only?"
$ws.Range("B9").Font.ThemeFont = 1
$ws.Range("B9").Font.Size = 20

$ws.Range("C9").Value = "yes"
$ws.Range("C9").Font.Name = "Arial"

# --- Row 10: leftover formatting left behind on the still-empty B10 -----
$ws.Range("B10").Font.ThemeFont = 1
$ws.Range("B10").Font.Size = 20

Write-Host "Added Q_7/Q_8 rows to SIQ sheet"
